$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("numeric")
$ws.Activate()

$ws.Range("E2").Value = "*"
$ws.Range("E3").Select()
